$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at 852 — this shifts the existing rows 852..912
# down to 853..913 (and grows the used range to A1:R913), matching every
# row of historical price data moving down by one slot.
$ws.Rows(852).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A852").Value = 4
$ws.Range("B852").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C852").Value = "Los Lagos"
$ws.Range("D852").Value = 45106
$ws.Range("E852").Value = 10
$ws.Range("F852").Value = 100112033
$ws.Range("G852").Value = "Lechuga"
$ws.Range("H852").Value = "Escarola"
$ws.Range("I852").Value = "Primera"
$ws.Range("J852").Value = 250
$ws.Range("K852").Value = 11000
$ws.Range("L852").Value = 11000
$ws.Range("M852").Value = 11000
$ws.Range("N852").Value = "$/caja 15 unidades"
$ws.Range("O852").Value = "Región de Coquimbo"
$ws.Range("P852").Value = 733
$ws.Range("Q852").Value = 15
$ws.Range("R852").Value = "Hortaliza"
